# seeds_UNCONFIRMED.xlsx cleanup:
#  - species code "jubu" / "juncus bufonius" renamed throughout the data to
#    "jute" / "juncus tenuis" (ID correction)
#  - the single unidentified "forb1" specimen (row 612, reference 6a) was
#    identified as Epilobium glaberrimum (epgl)
#  - sheet view scroll position / selection updated to where the editor
#    left off

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose G (species code) / H (species latin name) columns currently
# read "jubu" / "juncus bufonius" -> rename to "jute" / "juncus tenuis".
$jubuRows = @(130,141,142,158,168,171,175,178,182,197,200,211,222,224,230,
              278,291,298,301,310,313,316,321,327,332,334,337,345,359,388,
              397,553,556,563,569,572,591,618,621)

foreach ($r in $jubuRows) {
    $ws.Range("G$r").Value = "jute"
    $ws.Range("H$r").Value = "juncus tenuis"
}

# Row 612: the "forb1" placeholder was identified as Epilobium glaberrimum.
$ws.Range("H612").Value = "epilobium glaberrimum"
$ws.Range("G612").Value = "epgl"

# Update the saved view: drop the old frozen/scrolled top-left cell and move
# the remembered selection to L9.
$ws.Range("L9").Select()
